$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last existing data row (30) into two new rows, inserting
# copies so the new rows inherit the same cell styles used throughout the
# table (e.g. the style on column D and the boolean style on column I).
$ws.Rows.Item(30).Copy()
$ws.Rows.Item(31).Insert(-4121)
$excel.CutCopyMode = $false

$ws.Rows.Item(30).Copy()
$ws.Rows.Item(32).Insert(-4121)
$excel.CutCopyMode = $false

# Fill in row 32 (John Doe) first, then row 31 (Jane Smith), so that the
# shared-string table records "John Doe"/"john.doe@xyz.com" before
# "Jane Smith"/"jane.smith@xyz.com", matching how the entries were typed.
$ws.Cells.Item(32, 1).Value = 110031
$ws.Cells.Item(32, 2).Value = 9317596767
$ws.Cells.Item(32, 3).Value = "John Doe"
$ws.Cells.Item(32, 4).Value = "john.doe@xyz.com"
$ws.Cells.Item(32, 5).Value = 818876431

$ws.Cells.Item(31, 1).Value = 110030
$ws.Cells.Item(31, 2).Value = 9317596768
$ws.Cells.Item(31, 3).Value = "Jane Smith"
$ws.Cells.Item(31, 4).Value = "jane.smith@xyz.com"
$ws.Cells.Item(31, 5).Value = 818876432

$ws.Range("E28").Select()
